$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 65.45238095238095
$ws.Cells.Item(2, 8).Value = 17.52119047619048
$ws.Cells.Item(2, 9).Value = 167.8352380952381
$ws.Cells.Item(2, 10).Value = 7147.117142857142

$ws.Cells.Item(3, 7).Value = 90.5
$ws.Cells.Item(3, 8).Value = 24.57071428571428
$ws.Cells.Item(3, 9).Value = 186.4028571428571
$ws.Cells.Item(3, 10).Value = 9895.018571428573

$ws.Cells.Item(4, 7).Value = 96.78571428571429
$ws.Cells.Item(4, 8).Value = 28.86642857142857
$ws.Cells.Item(4, 9).Value = 197.6928571428571
$ws.Cells.Item(4, 10).Value = 12880.15142857143

$ws.Cells.Item(5, 7).Value = 68.0
$ws.Cells.Item(5, 8).Value = 17.10777777777778
$ws.Cells.Item(5, 9).Value = 124.8822222222222
$ws.Cells.Item(5, 10).Value = 6191.664444444445

$ws.Cells.Item(6, 7).Value = 133.7222222222222
$ws.Cells.Item(6, 8).Value = 33.53666666666666
$ws.Cells.Item(6, 9).Value = 128.4011111111111
$ws.Cells.Item(6, 10).Value = 13032.63444444444

$ws.Cells.Item(7, 7).Value = 111.3125
$ws.Cells.Item(7, 8).Value = 30.713125
$ws.Cells.Item(7, 9).Value = 155.031875
$ws.Cells.Item(7, 10).Value = 13145.564375

$ws.Cells.Item(8, 7).Value = 134.7142857142857
$ws.Cells.Item(8, 8).Value = 38.19928571428571
$ws.Cells.Item(8, 9).Value = 161.6435714285714
$ws.Cells.Item(8, 10).Value = 16368.36285714286

$ws.Cells.Item(9, 7).Value = 278.625
$ws.Cells.Item(9, 8).Value = 66.910625
$ws.Cells.Item(9, 9).Value = 117.86875
$ws.Cells.Item(9, 10).Value = 24165.876875

$ws.Cells.Item(10, 7).Value = 191.3125
$ws.Cells.Item(10, 8).Value = 52.52125
$ws.Cells.Item(10, 9).Value = 154.361875
$ws.Cells.Item(10, 10).Value = 20433.32125

$ws.Cells.Item(11, 7).Value = 165.9375
$ws.Cells.Item(11, 8).Value = 43.170625
$ws.Cells.Item(11, 9).Value = 125.0475
$ws.Cells.Item(11, 10).Value = 15489.34625

$ws.Cells.Item(12, 7).Value = 249.0
$ws.Cells.Item(12, 8).Value = 58.83428571428572
$ws.Cells.Item(12, 9).Value = 84.25142857142858
$ws.Cells.Item(12, 10).Value = 18264.09428571429

$ws.Cells.Item(13, 7).Value = 198.5625
$ws.Cells.Item(13, 8).Value = 52.12125
$ws.Cells.Item(13, 9).Value = 125.6175
$ws.Cells.Item(13, 10).Value = 19100.43375

$ws.Cells.Item(14, 7).Value = 106.5
$ws.Cells.Item(14, 8).Value = 26.948125
$ws.Cells.Item(14, 9).Value = 121.74625
$ws.Cells.Item(14, 10).Value = 9708.01125

$ws.Cells.Item(15, 7).Value = 135.9166666666667
$ws.Cells.Item(15, 8).Value = 27.24541666666667
$ws.Cells.Item(15, 9).Value = 63.51666666666667
$ws.Cells.Item(15, 10).Value = 7768.971666666666

$ws.Cells.Item(16, 7).Value = 149.8
$ws.Cells.Item(16, 8).Value = 28.383
$ws.Cells.Item(16, 9).Value = 46.726
$ws.Cells.Item(16, 10).Value = 7865.425999999999

$ws.Cells.Item(17, 7).Value = 79.4375
$ws.Cells.Item(17, 8).Value = 20.95125
$ws.Cells.Item(17, 9).Value = 150.763125
$ws.Cells.Item(17, 10).Value = 8334.363125

$ws.Cells.Item(18, 7).Value = 144.4444444444445
$ws.Cells.Item(18, 8).Value = 38.96777777777778
$ws.Cells.Item(18, 9).Value = 125.4611111111111
$ws.Cells.Item(18, 10).Value = 15320.04833333333

$ws.Cells.Item(19, 7).Value = 140.1111111111111
$ws.Cells.Item(19, 8).Value = 34.17555555555555
$ws.Cells.Item(19, 9).Value = 108.1377777777778
$ws.Cells.Item(19, 10).Value = 11662.68888888889

